# Adds a new "Tiles" task group (Tile Class / Tileset Class / Tilemap Class)
# to the bottom of the task list on Sheet1, mirroring the existing groups.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values, written in reading order so the shared-string table grows
# in the same order as the source edit (Tiles, Tile Class, Tileset Class,
# Tilemap Class).
$ws.Range("A10").Value = "Tiles"
$ws.Range("B10").Value = "Tile Class"
$ws.Range("B11").Value = "Tileset Class"
$ws.Range("B12").Value = "Tilemap Class"

# Group label spans the three new rows, like the other task groups above it.
$ws.Range("A10:A12").Merge()
$ws.Range("A10:A12").HorizontalAlignment = -4131
$ws.Range("A10:A12").VerticalAlignment = -4108

# Match the formatting already used by the other "Subtask" column entries.
$ws.Range("B10:B12").HorizontalAlignment = -4131
$ws.Range("B10:B12").VerticalAlignment = -4108

# Scroll the view down a bit and leave the last-entered cell selected,
# matching where the author ended up after adding the rows.
$excel.ActiveWindow.ScrollRow = 4
$ws.Range("B12").Select()
